$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the "Main Publication Output - PhD Researchers (ORCID)" heading
#    paragraph by scanning (robust to any paragraph-index drift).
# ---------------------------------------------------------------------------
$headingIdx = 0
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like 'Main Publication Output - PhD Researchers (ORCID)*') {
        $headingIdx = $i
        break
    }
}

if ($headingIdx -eq 0) {
    Write-Output "ERROR: heading paragraph not found"
} else {
    # -----------------------------------------------------------------
    # 2) Update the heading text (stays bold).
    # -----------------------------------------------------------------
    $headingRange = $d.Paragraphs.Item($headingIdx).Range
    $headingRange.Find.Execute(
        'Main Publication Output - PhD Researchers (ORCID)', $true, $false, $false, $false, $false,
        $true, 1, $false, 'Main Publication Output (ORCID 2024-2025)', 2) | Out-Null

    # -----------------------------------------------------------------
    # 3) Update the intro paragraph text (plain, not bold).
    # -----------------------------------------------------------------
    $introIdx = $headingIdx + 1
    $introRange = $d.Paragraphs.Item($introIdx).Range
    $introRange.Find.Execute(
        "The following publications are registered in ORCID for the project's PhD researchers:",
        $true, $false, $false, $false, $false, $true, 1, $false,
        'The following publications are registered in ORCID for the project researchers (2024-2025):', 2) | Out-Null

    # -----------------------------------------------------------------
    # 4) Insert the new Joerg Osterrieder block (name line + 17 numbered
    #    publications) right after the intro paragraph, before the
    #    "Lennart John Baals" sub-heading.
    # -----------------------------------------------------------------
    $newParas = @(
        'Joerg Osterrieder (ORCID: 0000-0003-0189-8636):',
        '1. "How can artificial intelligence help customer intelligence for credit portfolio management? A systematic literature review". International Journal of Information Management Data Insights. DOI: 10.1016/j.jjimei.2024.100234',
        '2. "Stylized facts of metaverse non-fungible tokens". Physica A: Statistical Mechanics and its Applications. DOI: 10.1016/j.physa.2024.130103',
        '3. "Leveraging network topology for credit risk assessment in P2P lending". Expert Systems with Applications. DOI: 10.1016/j.eswa.2024.124100 (with Baals, Liu)',
        '4. "Network centrality and credit risk: A comprehensive analysis of peer-to-peer lending dynamics". Finance Research Letters. DOI: 10.1016/j.frl.2024.105308 (with Baals, Liu)',
        '5. "Towards a new PhD Curriculum for Digital Finance". Open Research Europe. DOI: 10.12688/openreseurope.16513.1 (with Liu)',
        '6. "Visual XAI tool". Zenodo. DOI: 10.5281/zenodo.10934115',
        '7. "A discussion paper for possible approaches to building a statistically valid backtesting framework". SSRN. DOI: 10.2139/ssrn.4893677',
        '8. "Enhancing Security in Blockchain Networks: Anomalies, Frauds, and Advanced Detection Techniques". arXiv. DOI: 10.48550/arxiv.2402.11231',
        '9. "Ethical Artificial Intelligence, Fintech and Data Protection: A Path Forward for Training in Europe". SSRN. DOI: 10.2139/ssrn.4885037',
        '10. "Forecasting Commercial Customers Credit Risk Through Early Warning Signals Data". SSRN. DOI: 10.2139/ssrn.4754568',
        '11. "How can Consumers Without Credit History Benefit from Information Processing and Machine Learning Tools by Financial Institutions?". SSRN. DOI: 10.2139/ssrn.4730445',
        '12. "Hypothesizing Multimodal Influence: Assessing the Impact of Textual and Non-Textual Data on Financial Instrument Pricing Using NLP and Generative AI". SSRN. DOI: 10.2139/ssrn.4698153',
        '13. "Integrating Early Warning Systems with Customer Segmentation". SSRN. DOI: 10.2139/ssrn.4779632',
        '14. "Integration of Early Warning Systems and Customer Segmentation Methods in the Financial Industry - A Systematic Literature Review". SSRN. DOI: 10.2139/ssrn.4730479',
        '15. "Metaverse Non Fungible Tokens". SSRN. DOI: 10.2139/ssrn.4733153',
        '16. "Modeling Commodity Price Co-Movement: Building on Traditional Methods & Exploring Applications of Machine Learning Models". SSRN. DOI: 10.2139/ssrn.4730474',
        "17. ""Predicting Retail Customers' Distress: Early Warning Systems and Machine Learning Applications"". SSRN. DOI: 10.2139/ssrn.4730470"
    )

    $insertAfterPara = $d.Paragraphs.Item($introIdx)
    $insertRange = $insertAfterPara.Range
    foreach ($txt in $newParas) {
        $insertRange.InsertParagraphAfter()
        $insertRange = $d.Paragraphs.Item($introIdx + 1).Range
        $insertRange.InsertAfter($txt)
        $introIdx = $introIdx + 1
        $insertRange = $d.Paragraphs.Item($introIdx).Range
    }

    # -----------------------------------------------------------------
    # 5) Recompute index of the "Lennart John Baals" sub-heading, append
    #    the "Co-author..." note, then drop its two detail paragraphs.
    # -----------------------------------------------------------------
    $baalsIdx = 0
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like 'Lennart John Baals (ORCID: 0000-0002-7737-9675):*') {
            $baalsIdx = $i
            break
        }
    }

    $baalsRange = $d.Paragraphs.Item($baalsIdx).Range
    $baalsRange.End = $baalsRange.End - 1
    $baalsRange.InsertAfter(' Co-author on publications 3, 4 above.')

    # The next two paragraphs are the old Baals detail lines - delete them.
    $detail1 = $d.Paragraphs.Item($baalsIdx + 1)
    $detail2 = $d.Paragraphs.Item($baalsIdx + 2)
    $deleteRange = $d.Range($detail1.Range.Start, $detail2.Range.End)
    $deleteRange.Delete()

    # -----------------------------------------------------------------
    # 6) Recompute index of the "Yiting Liu" sub-heading, append the
    #    "Co-author..." note, then drop its four detail paragraphs.
    # -----------------------------------------------------------------
    $liuIdx = 0
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -like 'Yiting Liu (ORCID: 0009-0006-9554-8205):*') {
            $liuIdx = $i
            break
        }
    }

    $liuRange = $d.Paragraphs.Item($liuIdx).Range
    $liuRange.End = $liuRange.End - 1
    $liuRange.InsertAfter(' Co-author on publications 3, 4, 5 above.')

    # The next four paragraphs are the old Liu detail lines - delete them.
    $ldetail1 = $d.Paragraphs.Item($liuIdx + 1)
    $ldetail4 = $d.Paragraphs.Item($liuIdx + 4)
    $ldeleteRange = $d.Range($ldetail1.Range.Start, $ldetail4.Range.End)
    $ldeleteRange.Delete()

    Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
}
